# Update TPM-derived NATMI ligand-receptor metrics (Ifnk-Ifnar2) with
# newly recomputed values, per commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.062802
$ws.Cells.Item(2, 8).Value = 0.188406
$ws.Cells.Item(2, 9).Value = 0.0482629354765083
$ws.Cells.Item(2, 10).Value = 0.0482629354765083
$ws.Cells.Item(2, 13).Value = 13.11310966666667
$ws.Cells.Item(2, 14).Value = 39.339329
$ws.Cells.Item(2, 15).Value = 0.2602886552498481
$ws.Cells.Item(2, 16).Value = 0.2602886552498481
$ws.Cells.Item(2, 17).Value = 0.823529513286
$ws.Cells.Item(2, 18).Value = 7.411765619573999
$ws.Cells.Item(2, 19).Value = 0.01256229457359053
$ws.Cells.Item(2, 20).Value = 0.01256229457359053
$ws.Cells.Item(3, 7).Value = 0.062802
$ws.Cells.Item(3, 8).Value = 0.188406
$ws.Cells.Item(3, 9).Value = 0.0482629354765083
$ws.Cells.Item(3, 10).Value = 0.0482629354765083
$ws.Cells.Item(3, 15).Value = 0.5256625072946081
$ws.Cells.Item(3, 16).Value = 0.5256625072946082
$ws.Cells.Item(3, 17).Value = 1.663148124414
$ws.Cells.Item(3, 18).Value = 14.968333119726
$ws.Cells.Item(3, 19).Value = 0.02537001567197924
$ws.Cells.Item(3, 20).Value = 0.02537001567197925
$ws.Cells.Item(4, 7).Value = 0.062802
$ws.Cells.Item(4, 8).Value = 0.188406
$ws.Cells.Item(4, 9).Value = 0.0482629354765083
$ws.Cells.Item(4, 10).Value = 0.0482629354765083
$ws.Cells.Item(4, 13).Value = 10.783589
$ws.Cells.Item(4, 14).Value = 32.350767
$ws.Cells.Item(4, 15).Value = 0.2140488374555438
$ws.Cells.Item(4, 16).Value = 0.2140488374555438
$ws.Cells.Item(4, 17).Value = 0.6772309563779999
$ws.Cells.Item(4, 18).Value = 6.095078607401999
$ws.Cells.Item(4, 19).Value = 0.01033062523093852
$ws.Cells.Item(4, 20).Value = 0.01033062523093852
$ws.Cells.Item(5, 7).Value = 0.5410386666666667
$ws.Cells.Item(5, 9).Value = 0.415784756212054
$ws.Cells.Item(5, 10).Value = 0.415784756212054
$ws.Cells.Item(5, 13).Value = 13.11310966666667
$ws.Cells.Item(5, 14).Value = 39.339329
$ws.Cells.Item(5, 15).Value = 0.2602886552498481
$ws.Cells.Item(5, 16).Value = 0.2602886552498481
$ws.Cells.Item(5, 17).Value = 7.094699369907111
$ws.Cells.Item(5, 18).Value = 63.852294329164
$ws.Cells.Item(5, 19).Value = 0.1082240550678215
$ws.Cells.Item(5, 20).Value = 0.1082240550678215
$ws.Cells.Item(6, 7).Value = 0.5410386666666667
$ws.Cells.Item(6, 9).Value = 0.415784756212054
$ws.Cells.Item(6, 10).Value = 0.415784756212054
$ws.Cells.Item(6, 15).Value = 0.5256625072946081
$ws.Cells.Item(6, 16).Value = 0.5256625072946082
$ws.Cells.Item(6, 19).Value = 0.2185624574453057
$ws.Cells.Item(6, 20).Value = 0.2185624574453057
$ws.Cells.Item(7, 7).Value = 0.5410386666666667
$ws.Cells.Item(7, 9).Value = 0.415784756212054
$ws.Cells.Item(7, 10).Value = 0.415784756212054
$ws.Cells.Item(7, 13).Value = 10.783589
$ws.Cells.Item(7, 14).Value = 32.350767
$ws.Cells.Item(7, 15).Value = 0.2140488374555438
$ws.Cells.Item(7, 16).Value = 0.2140488374555438
$ws.Cells.Item(7, 17).Value = 5.834338614441333
$ws.Cells.Item(7, 18).Value = 52.50904752997199
$ws.Cells.Item(7, 19).Value = 0.08899824369892685
$ws.Cells.Item(7, 20).Value = 0.08899824369892687
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.6974063333333334
$ws.Cells.Item(8, 8).Value = 2.092219
$ws.Cells.Item(8, 9).Value = 0.5359523083114377
$ws.Cells.Item(8, 10).Value = 0.5359523083114377
$ws.Cells.Item(8, 13).Value = 13.11310966666667
$ws.Cells.Item(8, 14).Value = 39.339329
$ws.Cells.Item(8, 15).Value = 0.2602886552498481
$ws.Cells.Item(8, 16).Value = 0.2602886552498481
$ws.Cells.Item(8, 17).Value = 9.14516573122789
$ws.Cells.Item(8, 18).Value = 82.306491581051
$ws.Cells.Item(8, 19).Value = 0.1395023056084361
$ws.Cells.Item(8, 20).Value = 0.1395023056084361
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.6974063333333334
$ws.Cells.Item(9, 8).Value = 2.092219
$ws.Cells.Item(9, 9).Value = 0.5359523083114377
$ws.Cells.Item(9, 10).Value = 0.5359523083114377
$ws.Cells.Item(9, 15).Value = 0.5256625072946081
$ws.Cells.Item(9, 16).Value = 0.5256625072946082
$ws.Cells.Item(9, 17).Value = 18.468998363711
$ws.Cells.Item(9, 18).Value = 166.220985273399
$ws.Cells.Item(9, 19).Value = 0.2817300341773232
$ws.Cells.Item(9, 20).Value = 0.2817300341773232
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.6974063333333334
$ws.Cells.Item(10, 8).Value = 2.092219
$ws.Cells.Item(10, 9).Value = 0.5359523083114377
$ws.Cells.Item(10, 10).Value = 0.5359523083114377
$ws.Cells.Item(10, 13).Value = 10.783589
$ws.Cells.Item(10, 14).Value = 32.350767
$ws.Cells.Item(10, 15).Value = 0.2140488374555438
$ws.Cells.Item(10, 16).Value = 0.2140488374555438
$ws.Cells.Item(10, 17).Value = 7.520543264663666
$ws.Cells.Item(10, 18).Value = 67.684889381973
$ws.Cells.Item(10, 19).Value = 0.1147199685256784
$ws.Cells.Item(10, 20).Value = 0.1147199685256784